$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy formatting (bold font, border, centered alignment) from the existing header cell H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill in the data values for the new columns I (I0) and J (IF), rows 2-46
$ifValues = @{
    2 = @(9, 10)
    3 = @(6, 6)
    4 = @(4, 4)
    5 = @(8, 8)
    6 = @(5, 5)
    7 = @(9, 9)
    8 = @(6, 6)
    9 = @(6, 6)
    10 = @(6, 6)
    11 = @(3, 3)
    12 = @(5, 5)
    13 = @(8, 8)
    14 = @(7, 7)
    15 = @(7, 8)
    16 = @(7, 7)
    17 = @(7, 7)
    18 = @(6, 6)
    19 = @(7, 7)
    20 = @(7, 7)
    21 = @(7, 7)
    22 = @(8, 8)
    23 = @(6, 6)
    24 = @(6, 6)
    25 = @(8, 8)
    26 = @(7, 7)
    27 = @(6, 7)
    28 = @(12, 13)
    29 = @(8, 8)
    30 = @(7, 7)
    31 = @(8, 8)
    32 = @(7, 7)
    33 = @(5, 5)
    34 = @(7, 8)
    35 = @(7, 7)
    36 = @(8, 8)
    37 = @(6, 7)
    38 = @(7, 7)
    39 = @(1, 3)
    40 = @(1, 2)
    41 = @(1, 5)
    42 = @(1, 6)
    43 = @(1, 6)
    44 = @(1, 5)
    45 = @(1, 3)
    46 = @(4, 5)
}

foreach ($r in $ifValues.Keys) {
    $pair = $ifValues[$r]
    $ws.Cells.Item($r, 9).Value = $pair[0]
    $ws.Cells.Item($r, 10).Value = $pair[1]
}
